$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NroSiniestro value in F2 (trailing space preserved, keep as text)
$ws.Range("F2").Value = "'0420194406719 "

# Update the active selection to H6 as recorded in the saved view state
$ws.Range("H6").Select()
